$d = $word.ActiveDocument
$r = $null

# --- LogBack (paragraph 17) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>LogBack</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(17).Range
$r.InsertXML($snippet)

# --- Izmeneniya v konfig servere (paragraph 16) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Изменения в </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr><w:t>конфиг</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> сервере</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(16).Range
$r.InsertXML($snippet)

# --- Configure (paragraph 15) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="en-US"/></w:rPr><w:t>Configure</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(15).Range
$r.InsertXML($snippet)

# --- Monitoring (paragraph 14) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Мониторинг</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(14).Range
$r.InsertXML($snippet)

# --- Hystrix - turbina (paragraph 13) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hystrix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>турбина</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(13).Range
$r.InsertXML($snippet)

# --- Hystrix - metrics stream (paragraph 12) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hystrix</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>метрикс</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>стрим</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(12).Range
$r.InsertXML($snippet)

# --- Hystrix (paragraph 11) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Hystrix</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(11).Range
$r.InsertXML($snippet)

# --- Eureka - izmenenie nazvaniy (paragraph 9) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="en-US"/></w:rPr><w:t>Eureka</w:t></w:r><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> – изменение названий</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(9).Range
$r.InsertXML($snippet)

# --- Eureka (paragraph 8) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/><w:lang w:val="en-US"/></w:rPr><w:t>Eureka</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(8).Range
$r.InsertXML($snippet)

# --- Logirovanie + 4 new (paragraph 6) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Логирование</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Eureka</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="en-US"/></w:rPr><w:t>Zuul</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Config</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Monitoring config</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(6).Range
$r.InsertXML($snippet)

# --- Ochered (rabbit) (paragraph 5) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>Очередь (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>раббит</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/><w:lang w:val="ru-RU"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(5).Range
$r.InsertXML($snippet)

# --- Title (paragraph 1) ---
$snippet = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">План по переходу от монолита к </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>микросервисам</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r = $d.Paragraphs(1).Range
$r.InsertXML($snippet)

